# Edit script: rename stat sheets to human-friendly names, bump the
# "age-days" value in column E for every player row by one day, and
# fix the mis-placed "Playing Time" header/merge on the StandardStats
# and PlayingTime sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Rename the worksheets (Matches and Possession keep their names)
# ---------------------------------------------------------------
$renames = @{
    "StandardStats"    = "Standard Stats"
    "ShootingStats"     = "Shooting Stats"
    "PassingStats"      = "Passing Stats"
    "PassTypes"         = "Pass Types"
    "GoalShotCreation"  = "Goal & Shot Creation"
    "DefensiveActions"  = "Defensive Actions"
    "PlayingTime"       = "Playing Time"
    "MiscStats"         = "Miscellaneous Stats"
}

# Keep an ordered list of the *original* sheet names so we can still find
# each sheet by its old name while we rename them one at a time.
$statSheetOrder = @("StandardStats","ShootingStats","PassingStats","PassTypes","GoalShotCreation","DefensiveActions","Possession","PlayingTime","MiscStats")

foreach ($oldName in $statSheetOrder) {
    if ($renames.ContainsKey($oldName)) {
        $ws = $wb.Worksheets.Item($oldName)
        $ws.Name = $renames[$oldName]
    }
}

# ---------------------------------------------------------------
# 2. Bump every "NN-NNN" age value in column E (rows 4-31) by 1 day
#    across all of the per-player stat sheets.
# ---------------------------------------------------------------
$statSheetNames = @("Standard Stats","Shooting Stats","Passing Stats","Pass Types","Goal & Shot Creation","Defensive Actions","Possession","Playing Time","Miscellaneous Stats")

foreach ($sheetName in $statSheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($row = 4; $row -le 31; $row++) {
        $cell = $ws.Range("E$row")
        $val = $cell.Value2
        if ($val -match '^(\d+)-(\d{3})$') {
            $age = $Matches[1]
            $days = [int]$Matches[2] + 1
            $cell.Value = "{0}-{1:D3}" -f $age, $days
        }
    }
}

# ---------------------------------------------------------------
# 3. Fix the header: "Playing Time" label belongs in F1 (merged
#    F1:I1), not G1 (merged G1:I1), on the Standard Stats and
#    Playing Time sheets.
# ---------------------------------------------------------------
$headerFixSheets = @("Standard Stats","Playing Time")
foreach ($sheetName in $headerFixSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("G1:I1").UnMerge()
    $ws.Range("F1").Value = "Playing Time"
    $ws.Range("G1").Value = ""
    $ws.Range("F1:I1").Merge()
}
